$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (CellRef, NewValue) to apply.
# Generated from the target OOXML diff; values are written via
# Range.Value so Excel COM interop stores them as plain numbers.
$sheetEdits = @{}

$sheetEdits["ALC"] = @(
    @{ Ref = "H17"; Value = 1400.6957 }
    @{ Ref = "J17"; Value = 1400.6957 }
    @{ Ref = "L17"; Value = 4202.0871 }
    @{ Ref = "N17"; Value = -4538.0871 }
    @{ Ref = "H18"; Value = 1936.3684 }
    @{ Ref = "I18"; Value = 1699 }
    @{ Ref = "K18"; Value = 1699 }
    @{ Ref = "M18"; Value = -1415 }
    @{ Ref = "H74"; Value = 9758.333000000001 }
    @{ Ref = "I74"; Value = 9683.333000000001 }
    @{ Ref = "K74"; Value = 9683.333000000001 }
    @{ Ref = "M74"; Value = -8747.333000000001 }
    @{ Ref = "H77"; Value = 9758.333000000001 }
    @{ Ref = "I77"; Value = 9683.333000000001 }
    @{ Ref = "K77"; Value = 48416.665 }
    @{ Ref = "M77"; Value = -43736.665 }
    @{ Ref = "H98"; Value = 4754.3335 }
    @{ Ref = "I98"; Value = 4098.625 }
    @{ Ref = "K98"; Value = 4098.625 }
    @{ Ref = "M98"; Value = -2600.625 }
    @{ Ref = "H116"; Value = 4581.5 }
    @{ Ref = "I116"; Value = 4497.8 }
    @{ Ref = "K116"; Value = 4497.8 }
    @{ Ref = "M116"; Value = -1055.8 }
    @{ Ref = "H122"; Value = 4754.3335 }
    @{ Ref = "I122"; Value = 4098.625 }
    @{ Ref = "K122"; Value = 12295.875 }
    @{ Ref = "M122"; Value = -9845.875 }
    @{ Ref = "H125"; Value = 982 }
    @{ Ref = "J125"; Value = 998.6667 }
    @{ Ref = "L125"; Value = 8988.0003 }
    @{ Ref = "N125"; Value = -13908.0003 }
    @{ Ref = "H132"; Value = 1321.7333 }
    @{ Ref = "I132"; Value = 1273.3214 }
    @{ Ref = "K132"; Value = 3819.9642 }
    @{ Ref = "M132"; Value = -1289.9642 }
    @{ Ref = "H135"; Value = 1860.8334 }
    @{ Ref = "I135"; Value = 1626.5714 }
    @{ Ref = "J135"; Value = 2188.8 }
    @{ Ref = "K135"; Value = 14639.1426 }
    @{ Ref = "L135"; Value = 19699.2 }
    @{ Ref = "M135"; Value = -12104.1426 }
    @{ Ref = "N135"; Value = -24769.2 }
    @{ Ref = "H137"; Value = 1588.1052 }
    @{ Ref = "I137"; Value = 1604.3529 }
    @{ Ref = "J137"; Value = 1450 }
    @{ Ref = "K137"; Value = 4813.0587 }
    @{ Ref = "L137"; Value = 4350 }
    @{ Ref = "M137"; Value = -2263.0587 }
    @{ Ref = "N137"; Value = -9450 }
    @{ Ref = "H138"; Value = 7688.9185 }
    @{ Ref = "I138"; Value = 5086.391 }
    @{ Ref = "K138"; Value = 15259.173 }
    @{ Ref = "M138"; Value = -10119.173 }
    @{ Ref = "H141"; Value = 3025.1765 }
    @{ Ref = "I141"; Value = 2276.75 }
    @{ Ref = "J141"; Value = 15000 }
    @{ Ref = "K141"; Value = 6830.25 }
    @{ Ref = "L141"; Value = 45000 }
    @{ Ref = "M141"; Value = -1650.25 }
    @{ Ref = "N141"; Value = -55360 }
)

$sheetEdits["ARM"] = @(
    @{ Ref = "H2"; Value = 3362 }
    @{ Ref = "I2"; Value = 2209 }
    @{ Ref = "K2"; Value = 2209 }
    @{ Ref = "M2"; Value = -2096 }
    @{ Ref = "H74"; Value = 1635.3077 }
    @{ Ref = "I74"; Value = 1094.875 }
    @{ Ref = "K74"; Value = 1094.875 }
    @{ Ref = "M74"; Value = -220.875 }
    @{ Ref = "H77"; Value = 1635.3077 }
    @{ Ref = "I77"; Value = 1094.875 }
    @{ Ref = "K77"; Value = 5474.375 }
    @{ Ref = "M77"; Value = -1106.375 }
    @{ Ref = "H92"; Value = 65000 }
    @{ Ref = "J92"; Value = 65000 }
    @{ Ref = "L92"; Value = 65000 }
    @{ Ref = "N92"; Value = -69992 }
    @{ Ref = "H116"; Value = 3362 }
    @{ Ref = "I116"; Value = 2209 }
    @{ Ref = "K116"; Value = 2209 }
    @{ Ref = "M116"; Value = 85 }
    @{ Ref = "H132"; Value = 2107.9092 }
    @{ Ref = "I132"; Value = 1918.7 }
    @{ Ref = "J132"; Value = 4000 }
    @{ Ref = "K132"; Value = 5756.1 }
    @{ Ref = "L132"; Value = 12000 }
    @{ Ref = "M132"; Value = -3226.1 }
    @{ Ref = "N132"; Value = -17060 }
)

$sheetEdits["BSM"] = @(
    @{ Ref = "H3"; Value = 3362 }
    @{ Ref = "I3"; Value = 2209 }
    @{ Ref = "K3"; Value = 2209 }
    @{ Ref = "M3"; Value = -2095 }
    @{ Ref = "H134"; Value = 2391.9443 }
    @{ Ref = "I134"; Value = 2237.0667 }
    @{ Ref = "K134"; Value = 6711.2001 }
    @{ Ref = "M134"; Value = -4176.2001 }
)

$sheetEdits["CRP"] = @(
    @{ Ref = "H22"; Value = 195.8 }
    @{ Ref = "I22"; Value = 195.8 }
    @{ Ref = "K22"; Value = 195.8 }
    @{ Ref = "M22"; Value = 154.2 }
    @{ Ref = "H58"; Value = 2805.9048 }
    @{ Ref = "I58"; Value = 1490.2 }
    @{ Ref = "K58"; Value = 1490.2 }
    @{ Ref = "M58"; Value = -1287.2 }
    @{ Ref = "H132"; Value = 2696.3333 }
    @{ Ref = "I132"; Value = 2844.2 }
    @{ Ref = "K132"; Value = 8532.599999999999 }
    @{ Ref = "M132"; Value = -6002.599999999999 }
    @{ Ref = "H134"; Value = 2072.1562 }
    @{ Ref = "I134"; Value = 1749.04 }
    @{ Ref = "J134"; Value = 3226.1428 }
    @{ Ref = "K134"; Value = 5247.12 }
    @{ Ref = "L134"; Value = 9678.428400000001 }
    @{ Ref = "M134"; Value = -2712.12 }
    @{ Ref = "N134"; Value = -14748.4284 }
    @{ Ref = "H136"; Value = 2805.9048 }
    @{ Ref = "I136"; Value = 1490.2 }
    @{ Ref = "K136"; Value = 4470.6 }
    @{ Ref = "M136"; Value = -1920.6 }
)

$sheetEdits["CUL"] = @(
    @{ Ref = "H44"; Value = 1264.2142 }
    @{ Ref = "I44"; Value = 275 }
    @{ Ref = "J44"; Value = 1429.0834 }
    @{ Ref = "K44"; Value = 825 }
    @{ Ref = "L44"; Value = 4287.2502 }
    @{ Ref = "M44"; Value = -427 }
    @{ Ref = "N44"; Value = -5083.2502 }
    @{ Ref = "I99"; Value = 3539 }
    @{ Ref = "J99"; Value = 609 }
    @{ Ref = "K99"; Value = 10617 }
    @{ Ref = "L99"; Value = 1827 }
    @{ Ref = "M99"; Value = -8371 }
    @{ Ref = "N99"; Value = -6319 }
    @{ Ref = "H103"; Value = 184.8 }
    @{ Ref = "J103"; Value = 99 }
    @{ Ref = "L103"; Value = 297 }
    @{ Ref = "N103"; Value = -2055 }
    @{ Ref = "H131"; Value = 1115.7894 }
    @{ Ref = "I131"; Value = 530.5 }
    @{ Ref = "J131"; Value = 1385.9231 }
    @{ Ref = "K131"; Value = 1591.5 }
    @{ Ref = "L131"; Value = 4157.7693 }
    @{ Ref = "M131"; Value = 3448.5 }
    @{ Ref = "N131"; Value = -14237.7693 }
    @{ Ref = "H137"; Value = 5636.6665 }
    @{ Ref = "I137"; Value = 4947.6 }
    @{ Ref = "J137"; Value = 6498 }
    @{ Ref = "K137"; Value = 14842.8 }
    @{ Ref = "L137"; Value = 19494 }
    @{ Ref = "M137"; Value = -9742.800000000001 }
    @{ Ref = "N137"; Value = -29694 }
    @{ Ref = "H140"; Value = 1602.6666 }
    @{ Ref = "I140"; Value = 1602.6666 }
    @{ Ref = "K140"; Value = 4807.9998 }
    @{ Ref = "M140"; Value = 372.0002000000004 }
)

$sheetEdits["GSM"] = @(
    @{ Ref = "H70"; Value = 7784 }
    @{ Ref = "I70"; Value = 7443.5 }
    @{ Ref = "K70"; Value = 7443.5 }
    @{ Ref = "M70"; Value = -7173.5 }
    @{ Ref = "H73"; Value = 7784 }
    @{ Ref = "I73"; Value = 7443.5 }
    @{ Ref = "K73"; Value = 7443.5 }
    @{ Ref = "M73"; Value = -6507.5 }
    @{ Ref = "H122"; Value = 49931.477 }
    @{ Ref = "J122"; Value = 169204 }
    @{ Ref = "L122"; Value = 507612 }
    @{ Ref = "N122"; Value = -512512 }
    @{ Ref = "H126"; Value = 4586.3335 }
    @{ Ref = "I126"; Value = 4300.8 }
    @{ Ref = "K126"; Value = 12902.4 }
    @{ Ref = "M126"; Value = -10432.4 }
    @{ Ref = "H132"; Value = 2302.3416 }
    @{ Ref = "J132"; Value = 3274.5 }
    @{ Ref = "L132"; Value = 9823.5 }
    @{ Ref = "N132"; Value = -14883.5 }
    @{ Ref = "H135"; Value = 90468 }
    @{ Ref = "J135"; Value = 90468 }
    @{ Ref = "L135"; Value = 90468 }
    @{ Ref = "N135"; Value = -100608 }
)

$sheetEdits["LTW"] = @(
    @{ Ref = "H16"; Value = 13271.143 }
    @{ Ref = "I16"; Value = 15199.6 }
    @{ Ref = "J16"; Value = 12199.777 }
    @{ Ref = "K16"; Value = 15199.6 }
    @{ Ref = "L16"; Value = 12199.777 }
    @{ Ref = "M16"; Value = -15029.6 }
    @{ Ref = "N16"; Value = -12539.777 }
    @{ Ref = "H40"; Value = 3044.182 }
    @{ Ref = "I40"; Value = 2717.5715 }
    @{ Ref = "J40"; Value = 3615.75 }
    @{ Ref = "K40"; Value = 2717.5715 }
    @{ Ref = "L40"; Value = 3615.75 }
    @{ Ref = "M40"; Value = -2581.5715 }
    @{ Ref = "N40"; Value = -3887.75 }
    @{ Ref = "H122"; Value = 7083.3335 }
    @{ Ref = "I122"; Value = 7083.3335 }
    @{ Ref = "K122"; Value = 21250.0005 }
    @{ Ref = "M122"; Value = -18800.0005 }
    @{ Ref = "H132"; Value = 5213 }
    @{ Ref = "I132"; Value = 4544.143 }
    @{ Ref = "K132"; Value = 13632.429 }
    @{ Ref = "M132"; Value = -11102.429 }
)

$sheetEdits["WVR"] = @(
    @{ Ref = "H62"; Value = 7587.8887 }
    @{ Ref = "I62"; Value = 5999.6665 }
    @{ Ref = "K62"; Value = 5999.6665 }
    @{ Ref = "M62"; Value = -5375.6665 }
    @{ Ref = "H65"; Value = 7587.8887 }
    @{ Ref = "I65"; Value = 5999.6665 }
    @{ Ref = "K65"; Value = 29998.3325 }
    @{ Ref = "M65"; Value = -26878.3325 }
    @{ Ref = "H81"; Value = 8832.833000000001 }
    @{ Ref = "I81"; Value = 5484.4287 }
    @{ Ref = "K81"; Value = 10968.8574 }
    @{ Ref = "M81"; Value = -9907.857400000001 }
    @{ Ref = "H84"; Value = 8832.833000000001 }
    @{ Ref = "I84"; Value = 5484.4287 }
    @{ Ref = "K84"; Value = 54844.287 }
    @{ Ref = "M84"; Value = -49540.287 }
    @{ Ref = "H100"; Value = 1421.25 }
    @{ Ref = "I100"; Value = 1504.7858 }
    @{ Ref = "J100"; Value = 1226.3334 }
    @{ Ref = "K100"; Value = 3009.5716 }
    @{ Ref = "L100"; Value = 2452.6668 }
    @{ Ref = "M100"; Value = -2468.5716 }
    @{ Ref = "N100"; Value = -3534.6668 }
    @{ Ref = "H132"; Value = 73472.234 }
    @{ Ref = "I132"; Value = 73472.234 }
    @{ Ref = "K132"; Value = 220416.702 }
    @{ Ref = "M132"; Value = -217886.702 }
    @{ Ref = "H136"; Value = 1676.2 }
    @{ Ref = "I136"; Value = 1659.1578 }
    @{ Ref = "K136"; Value = 4977.4734 }
    @{ Ref = "M136"; Value = -2427.4734 }
)

foreach ($sheetName in $sheetEdits.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($edit in $sheetEdits[$sheetName]) {
        $ws.Range($edit.Ref).Value = $edit.Value
    }
}

Write-Output "Applied $($sheetEdits.Values | ForEach-Object { $_.Count } | Measure-Object -Sum | Select-Object -ExpandProperty Sum) cell edits across $($sheetEdits.Keys.Count) sheets"